$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
#    Structure to produce:
#      <w:p>
#        <w:r/>
#        <w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>
#        <w:r><w:t>: Play Dr Fortuno, ...</w:t></w:r>
#      </w:p>
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$insertionPoint = $d.Range($titlePara.Range.End, $titlePara.Range.End)

$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"
$metaXml = "<w:p $wns><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Play Dr Fortuno, a circus-themed online slot game with special features and bonuses. Read our review and play for free.</w:t></w:r></w:p><w:p $wns></w:p>"
$insertionPoint.InsertXML($metaXml)

# InsertXML splices the *last* xml paragraph's runs into whatever
# paragraph used to follow the insertion point; since that trailing
# xml paragraph carries no runs at all, it is left behind as its own
# now-empty paragraph immediately after our new "Meta description"
# paragraph (rather than merging into "Game Overview and Design").
# Remove that leftover empty spacer paragraph.
$spacerPara = $d.Paragraphs.Item(2).Next()
$spacerPara.Range.Delete()

# ------------------------------------------------------------------
# 2. Remove the duplicated bold "Play Dr Fortuno Free Slot Game -
#    Review" paragraph that used to sit near the end of the document
#    (right before the closing italic paragraph) -- that heading now
#    only needs to appear once, at the very top of the document.
# ------------------------------------------------------------------
$dupPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Play Dr Fortuno Free Slot Game - Review" -and $i -ne 1) {
        $dupPara = $p
    }
}
if ($dupPara -ne $null) {
    $dupPara.Range.Delete()
}

# ------------------------------------------------------------------
# 3. Update the closing italic paragraph's text (image-prompt line).
#    Assign Range.Text directly (rather than Find/Replace) so smart
#    quotes / autocorrect don't mangle the straight quotes/apostrophes
#    required by the target text.
# ------------------------------------------------------------------
$oldText = "Play Dr Fortuno, a circus-themed online slot game with special features and bonuses. Read our review and play for free."
$newText = 'Create an eye-catching cartoon-style feature image for the game "Dr Fortuno". The image should include a happy Maya warrior wearing glasses, who embodies the theme of the game. The warrior should be in a celebratory pose, and the image should be bright and colorful to attract the attention of prospective players. Be sure to include the game''s title, "Dr Fortuno," in the image, as well as any other relevant design elements that reflect the game''s circus-inspired theme.'

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastText = $lastPara.Range.Text.TrimEnd([char]13, [char]7)
if ($lastText -eq $oldText) {
    $start = $lastPara.Range.Start
    $end = $lastPara.Range.End - 1
    $target = $d.Range($start, $end)
    $target.Text = $newText
}

Write-Output "Done. Paragraph count: $($d.Paragraphs.Count)"
